# The source data regenerated the stimulus/order table with updated
# distance and size codes:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31
# These codes appear (sometimes combined, e.g. "Face17_D64_S30") in the
# Condition, Filename_Left, Filename_Right, Distance and Size columns, as
# well as inside the fixation/face image filenames. Apply the substitution
# to every used cell on the sheet via Find/Replace. The four old tokens are
# disjoint from each other and from the new tokens, so running the four
# replacements in sequence over the whole used range correctly rewrites
# every occurrence exactly once.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

$used.Replace("D64", "D69") | Out-Null
$used.Replace("D80", "D86") | Out-Null
$used.Replace("D51", "D55") | Out-Null
$used.Replace("S30", "S31") | Out-Null
